# Weekly update for "Hortaliza, Vega Central Mapocho de Santiago - Zapallo":
# a new week of price data (2022-01-17) is inserted at the top of the
# Camote/Zapallo block (row 797), pushing the existing rows down by four
# and extending the used range from A1:R857 to A1:R861.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 797:857 down by 4 rows, opening up 4 blank rows at 797:800.
$ws.Rows.Item(797).Resize(4).Insert()

# --- Fill the 4 new rows (797:800) with the new week's data ---------------

# Columns shared by all four new rows.
$ws.Range("A797:A800").Value = 9
$ws.Range("B797:B800").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C797:C800").Value = "Metropolitana"
$ws.Range("D797:D800").Value = 44578
$ws.Range("E797:E800").Value = 13
$ws.Range("F797:F800").Value = 100112045
$ws.Range("G797:G800").Value = "Zapallo"
$ws.Range("H797:H800").Value = "Camote"
$ws.Range("N797:N800").Value = "$/kilo (volumen en unidades)"
$ws.Range("Q797:Q800").Value = 1
$ws.Range("R797:R800").Value = "Hortaliza"

# Row 797: 1a nueva(o) - Región Metropolitana
$ws.Range("I797").Value = "1a nueva(o)"
$ws.Range("J797").Value = 97
$ws.Range("K797").Value = 450
$ws.Range("L797").Value = 500
$ws.Range("M797").Value = 475
$ws.Range("O797").Value = "Región Metropolitana"
$ws.Range("P797").Value = 475

# Row 798: 1a nueva(o) - Región de O'Higgins
$ws.Range("I798").Value = "1a nueva(o)"
$ws.Range("J798").Value = 106
$ws.Range("K798").Value = 450
$ws.Range("L798").Value = 500
$ws.Range("M798").Value = 475
$ws.Range("O798").Value = "Región de O'Higgins"
$ws.Range("P798").Value = 475

# Row 799: 2a nueva(o) - Región Metropolitana
$ws.Range("I799").Value = "2a nueva(o)"
$ws.Range("J799").Value = 43
$ws.Range("K799").Value = 400
$ws.Range("L799").Value = 400
$ws.Range("M799").Value = 400
$ws.Range("O799").Value = "Región Metropolitana"
$ws.Range("P799").Value = 400

# Row 800: 2a nueva(o) - Región de O'Higgins
$ws.Range("I800").Value = "2a nueva(o)"
$ws.Range("J800").Value = 79
$ws.Range("K800").Value = 400
$ws.Range("L800").Value = 400
$ws.Range("M800").Value = 400
$ws.Range("O800").Value = "Región de O'Higgins"
$ws.Range("P800").Value = 400
